$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CustomerMappingDriver Class grading (rows 29-30) ---
# Row 30: "For no output" -> 0 scored (typed first, becomes the first new shared string)
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = "For no output"

# Row 29: "For incorrect logic to scan data" -> full marks scored (8 out of 8 available)
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = "For incorrect logic to scan data"

# --- Final section (row 37): compilation-error deduction comment ---
$ws.Range("E37").Value = -5
$ws.Range("F37").Value = "Compilation errors"

# Move the on-screen selection to match the author's final cursor position
$excel.Goto($ws.Range("F37"))
